$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 199-200, pushing the existing rows 199-296 down to 201-298.
$ws.Range("A199:A200").EntireRow.Insert()

# New row 199: Primera quality, week of 2022-02-17 (serial 44609)
$ws.Range("A199").Value = 3
$ws.Range("B199").Value = "Femacal de La Calera"
$ws.Range("C199").Value = "Coquimbo"
$ws.Range("D199").Value = 44609
$ws.Range("E199").Value = 5
$ws.Range("F199").Value = 100112043
$ws.Range("G199").Value = "Pepino ensalada"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 85
$ws.Range("K199").Value = 16000
$ws.Range("L199").Value = 17000
$ws.Range("M199").Value = 16529
$ws.Range("N199").Value = "$/caja 70 unidades"
$ws.Range("O199").Value = "Limache"
$ws.Range("P199").Value = 236
$ws.Range("Q199").Value = 70
$ws.Range("R199").Value = "Hortaliza"

# New row 200: Segunda quality, same week (serial 44609)
$ws.Range("A200").Value = 3
$ws.Range("B200").Value = "Femacal de La Calera"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44609
$ws.Range("E200").Value = 5
$ws.Range("F200").Value = 100112043
$ws.Range("G200").Value = "Pepino ensalada"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Segunda"
$ws.Range("J200").Value = 40
$ws.Range("K200").Value = 13000
$ws.Range("L200").Value = 13000
$ws.Range("M200").Value = 13000
$ws.Range("N200").Value = "$/caja 100 unidades"
$ws.Range("O200").Value = "Limache"
$ws.Range("P200").Value = 130
$ws.Range("Q200").Value = 100
$ws.Range("R200").Value = "Hortaliza"
